# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted before
# column N, pushing the old N/O/P columns (Late / heading / Outstanding) one
# slot to the right (-> O/P/Q). The new column is given roughly the same
# width as its neighbours. The workbook's active sheet/selection also moves
# from "Input" (A1:B2, cell B3 selected) to "Repayment schedule", with cell
# S13 selected there.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item(3)

# Insert a blank column before column N, shifting N:P -> O:Q.
$wsSchedule.Columns("N:N").Insert()

# Give the freshly inserted column the same width as column M/O (~10.71 chars).
$wsSchedule.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet, with S13 selected.
$wsSchedule.Range("S13").Select() | Out-Null
$wsSchedule.Activate() | Out-Null
